# Update betting-odds values in row 2 of Sheet1 to reflect the latest
# FlashScore data refresh (Jogos_da_Semana_FlashScore_2024-11-26.xlsx).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.7
$ws.Range("H2").Value = 3.6
$ws.Range("I2").Value = 4.75
$ws.Range("K2").Value = 2.2
$ws.Range("L2").Value = 5
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 11
$ws.Range("O2").Value = 1.3
$ws.Range("P2").Value = 3.4
$ws.Range("Q2").Value = 1.98
$ws.Range("R2").Value = 1.83
$ws.Range("S2").Value = 1.4
$ws.Range("T2").Value = 2.75
$ws.Range("U2").Value = 1.8
$ws.Range("V2").Value = 1.91
$ws.Range("W2").Value = 7
$ws.Range("X2").Value = 8
$ws.Range("Z2").Value = 13
$ws.Range("AB2").Value = 26
$ws.Range("AC2").Value = 11
$ws.Range("AF2").Value = 51
$ws.Range("AG2").Value = 251
$ws.Range("AH2").Value = 13
$ws.Range("AL2").Value = 41
$ws.Range("AP2").Value = 21
$ws.Range("AR2").Value = 51
$ws.Range("AS2").Value = 151
$ws.Range("AT2").Value = 2.75
$ws.Range("AU2").Value = 8.5
$ws.Range("AX2").Value = 26
$ws.Range("AY2").Value = 34
$ws.Range("BA2").Value = 101
$ws.Range("BB2").Value = 201
